$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) from "ShearF-HW10.xpc" to "ShearF"
$ws.Name = "ShearF"

# Add a new row 16 with data, mirroring the pattern of the existing rows
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.998451683503045
$ws.Range("D16").Value = 4.564337571077123
$ws.Range("E16").Value = 0.5187705011526014
$ws.Range("F16").Value = 1.998451683503045
$ws.Range("G16").Value = 1.188758491135423
$ws.Range("H16").Value = 0.2971667452826968
$ws.Range("I16").Value = 0.6188239908525918
$ws.Range("J16").Value = 4.564337571077123
$ws.Range("K16").Value = 2.541554036114862
$ws.Range("L16").Value = 2.270002859808954
$ws.Range("M16").Value = 1.531051497167247
